$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Forming the consolidated report: mark absent-day totals in column H
# and ensure previously-blank attendance cells are populated with 0.
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
